# Generate Report for Handoff
# Update status from "In Translation" to "Ready for handoff" and refresh the
# related timestamps on the Overview / zh-cn / de-de report sheets, widening
# the Status columns so the new, longer status text fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-13 02:46:50"

# --- zh-cn sheet ------------------------------------------------------------
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-13 02:46:43"

# --- de-de sheet ------------------------------------------------------------
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-13 02:46:50"

# --- Widen the Status columns to fit the new text --------------------------
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
